$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CalculationAmountAndPrice")
$ws2 = $wb.Worksheets.Item("CalculationLiquidity")

# --- Sheet1: CalculationAmountAndPrice ---

# Row2: relabel the trade token (18-decimal SSRT) and shrink the swap amount
$ws1.Range("C2").Value2 = "Snip 20 - SSRT ( 18 decimals)"
$ws1.Range("D2").Value2 = 100

# Fix the sell-side price formulas (J4, J5): quote price as 1/(F/D) instead of F/D
$ws1.Range("J4").Formula = "=1/ (F4/D4)"
$ws1.Range("J5").Formula = "=1 / (F5/D5)"

# Recalculate the dependency chain (A3:B7, F2:F6, J2:J6, etc.)
$ws1.Calculate()

# Widen column C to fit the new longer label
$ws1.Columns.Item(3).ColumnWidth = 26.16666666666667

# --- Sheet2: CalculationLiquidity ---
# Update the selection on sheet2 (leave sheet1 as the selected tab afterwards)
$ws2.Activate()
$ws2.Range("K5").Select()
$ws2.Columns.Item(13).ColumnWidth = 14.592447916666666

# --- Re-activate sheet1 (it was the selected tab originally), set zoom + selection ---
$ws1.Activate()
$excel.ActiveWindow.Zoom = 100
$ws1.Range("J6").Select()

$wb.Save()
